$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format while writing, so numeric-looking strings
# (e.g. "0.390", "14.80", "169.30", "1.00") keep their exact text, then
# restore the default "Normal" style so no stray formatting is introduced.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "37.723.02"
$ws.Range("E2").Value = "  +0.12%  "

$ws.Range("D3").Value = "2.076.27"
$ws.Range("E3").Value = "  -1.43%  "

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").Value = "233.51"
$ws.Range("E5").Value = "  -0.69%  "

$ws.Range("D6").Value = "0.623"
$ws.Range("E6").Value = "  -0.30%  "

$ws.Range("E7").Value = "  -0.13%  "

$ws.Range("D8").Value = "58.03"
$ws.Range("E8").Value = "  -0.11%  "

$ws.Range("D9").Value = "0.390"
$ws.Range("E9").Value = "  -0.28%  "

$ws.Range("E10").Value = "  +0.37%  "

$ws.Range("D11").Value = "0.105"
$ws.Range("E11").Value = "  +2.66%  "

$ws.Range("D12").Value = "2.382.69"
$ws.Range("E12").Value = "  -1.69%  "

$ws.Range("D13").Value = "14.80"
$ws.Range("E13").Value = "  +1.45%  "

$ws.Range("D14").Value = "20.85"
$ws.Range("E14").Value = "  -1.05%  "

$ws.Range("D15").Value = "0.772"
$ws.Range("E15").Value = "  -1.82%  "

$ws.Range("E16").Value = "  +1.23%  "

$ws.Range("D17").Value = "2.045.96"
$ws.Range("E17").Value = "  -2.67%  "

$ws.Range("D18").Value = "37.672.18"
$ws.Range("E18").Value = "  +0.03%  "

$ws.Range("E19").Value = "  -0.11%  "

$ws.Range("D20").Value = "71.04"
$ws.Range("E20").Value = "  +1.50%  "

$ws.Range("E21").Value = "  +1.28%  "

$ws.Range("D22").Value = "227.57"
$ws.Range("E22").Value = "  +0.12%  "

$ws.Range("E24").Value = "  -0.56%  "

$ws.Range("E25").Value = "  -1.23%  "

$ws.Range("D26").Value = "169.30"
$ws.Range("E26").Value = "  +0.78%  "

$ws.Range("D27").Value = "0.138"
$ws.Range("E27").Value = "  +3.76%  "

$ws.Range("D28").Value = "8.97"
$ws.Range("E28").Value = "  +0.23%  "

$ws.Range("D29").Value = "19.44"
$ws.Range("E29").Value = "  -0.04%  "

$ws.Range("E30").Value = "  -1.72%  "

$ws.Range("E31").Value = "  +2.25%  "

$ws.Range("E32").Value = "  +0.86%  "

$ws.Range("D33").Value = "0.0628"
$ws.Range("E33").Value = "  +1.19%  "

$ws.Range("D34").Value = "4.65"
$ws.Range("E34").Value = "  +1.46%  "

$ws.Range("E35").Value = "  -3.68%  "

$ws.Range("D36").Value = "1.82"
$ws.Range("E36").Value = "  +2.69%  "

$ws.Range("E37").Value = "  -2.97%  "

$ws.Range("D38").Value = "1.00"
$ws.Range("E38").Value = "  -0.11%  "

$ws.Range("E39").Value = "  -5.09%  "

$ws.Range("D40").Value = "0.0976"
$ws.Range("E40").Value = "  +1.49%  "

$ws.Range("D41").Value = "97.89"
$ws.Range("E41").Value = "  +0.43%  "

$ws.Range("B42").Value = "HuobiToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D42").Value = "2.87"
$ws.Range("E42").Value = "  -2.78%  "

$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").Value = "0.0215"
$ws.Range("E43").Value = "  +0.51%  "

$ws.Range("D44").Value = "1.449.29"
$ws.Range("E44").Value = "  -1.57%  "

$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").Value = "16.48"
$ws.Range("E45").Value = "  +6.08%  "

$ws.Range("B46").Value = "TrustWalletToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D46").Value = "1.16"
$ws.Range("E46").Value = "  -0.62%  "

$ws.Range("D47").Value = "4.26"
$ws.Range("E47").Value = "  +1.12%  "

$ws.Range("E48").Value = "  +1.15%  "

$ws.Range("E49").Value = "  +0.92%  "

$ws.Range("E50").Value = "  -0.51%  "

$ws.Range("D51").Value = "2.266.47"
$ws.Range("E51").Value = "  -1.65%  "

$ws.Range("D2:D51").Style = "Normal"
